# Regenerate the "K" column (column G) values in the save_data sheet.
# Mirrors an upstream data regen that switched the source of this column
# from "Strike#" to "K" and recalculated std/mean, writing refreshed s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2-18 (column G), replacing the old data.
$kValues = @{
    2  = 0
    3  = 1
    4  = 3
    5  = 0
    6  = 0
    7  = 1
    8  = 2
    9  = 0
    10 = 1
    11 = 1
    12 = 2
    13 = 2
    14 = 1
    15 = 2
    16 = 0
    17 = 2
    18 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
